$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H87").Value = 17316.89
$ws.Range("J87").Value = 17316.89
$ws.Range("L87").Value = 17316.89
$ws.Range("N87").Value = -19812.89
$ws.Range("H90").Value = 17316.89
$ws.Range("J90").Value = 17316.89
$ws.Range("L90").Value = 51950.67
$ws.Range("N90").Value = -64430.67
$ws.Range("H113").Value = 9175.549999999999
$ws.Range("I113").Value = 2945
$ws.Range("J113").Value = 14273.272
$ws.Range("K113").Value = 2945
$ws.Range("L113").Value = 14273.272
$ws.Range("M113").Value = 309
$ws.Range("N113").Value = -20781.272
$ws.Range("H116").Value = 5392.857
$ws.Range("I116").Value = 8266.666999999999
$ws.Range("J116").Value = 3237.5
$ws.Range("K116").Value = 8266.666999999999
$ws.Range("L116").Value = 3237.5
$ws.Range("M116").Value = -4824.666999999999
$ws.Range("N116").Value = -10121.5
$ws.Range("H121").Value = 2567.2222
$ws.Range("I121").Value = 798.75
$ws.Range("J121").Value = 3982
$ws.Range("K121").Value = 2396.25
$ws.Range("L121").Value = 11946
$ws.Range("M121").Value = -649.25
$ws.Range("N121").Value = -15440
$ws.Range("H129").Value = 723.8823
$ws.Range("J129").Value = 847.1667
$ws.Range("L129").Value = 2541.5001
$ws.Range("N129").Value = -12541.5001
$ws.Range("H132").Value = 1572.25
$ws.Range("I132").Value = 764.439
$ws.Range("J132").Value = 4583.1816
$ws.Range("K132").Value = 2293.317
$ws.Range("L132").Value = 13749.5448
$ws.Range("M132").Value = 236.683
$ws.Range("N132").Value = -18809.5448

$ws = $wb.Worksheets("ARM")
$ws.Range("H94").Value = 19999
$ws.Range("J94").Value = 19999
$ws.Range("L94").Value = 19999
$ws.Range("N94").Value = -21801
$ws.Range("H122").Value = 2176.35
$ws.Range("I122").Value = 1937.7142
$ws.Range("J122").Value = 2733.1667
$ws.Range("K122").Value = 5813.142599999999
$ws.Range("L122").Value = 8199.500100000001
$ws.Range("M122").Value = -3363.142599999999
$ws.Range("N122").Value = -13099.5001

$ws = $wb.Worksheets("BSM")
$ws.Range("H18").Value = 3900
$ws.Range("J18").Value = 3900
$ws.Range("L18").Value = 3900
$ws.Range("N18").Value = -4958
$ws.Range("H94").Value = 566.5
$ws.Range("I94").Value = 556.44446
$ws.Range("J94").Value = 596.6667
$ws.Range("K94").Value = 556.44446
$ws.Range("L94").Value = 596.6667
$ws.Range("M94").Value = -105.44446
$ws.Range("N94").Value = -1498.6667
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("K97").Value = 1000
$ws.Range("M97").Value = -9

$ws = $wb.Worksheets("CRP")
$ws.Range("H132").Value = 2840.8518
$ws.Range("I132").Value = 1723.6666
$ws.Range("J132").Value = 4237.3335
$ws.Range("K132").Value = 5170.9998
$ws.Range("L132").Value = 12712.0005
$ws.Range("M132").Value = -2640.9998
$ws.Range("N132").Value = -17772.0005
$ws.Range("H134").Value = 3967
$ws.Range("I134").Value = 5125.5
$ws.Range("J134").Value = 1650
$ws.Range("K134").Value = 15376.5
$ws.Range("L134").Value = 4950
$ws.Range("M134").Value = -12841.5
$ws.Range("N134").Value = -10020

$ws = $wb.Worksheets("CUL")
$ws.Range("H5").Value = 1254281.9
$ws.Range("I5").Value = 388.6
$ws.Range("K5").Value = 1165.8
$ws.Range("M5").Value = -1053.8
$ws.Range("H87").Value = 5000
$ws.Range("I87").Value = 5000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 15000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -13752
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 5000
$ws.Range("I90").Value = 5000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 45000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -38760
$ws.Range("N90").ClearContents()
$ws.Range("H113").Value = 1421.4
$ws.Range("I113").Value = 1400.1
$ws.Range("J113").Value = 1442.7
$ws.Range("K113").Value = 4200.299999999999
$ws.Range("L113").Value = 4328.1
$ws.Range("M113").Value = -2030.299999999999
$ws.Range("N113").Value = -8668.1
$ws.Range("H122").Value = 952.5217
$ws.Range("I122").Value = 492.16666
$ws.Range("J122").Value = 2609.8
$ws.Range("K122").Value = 4429.49994
$ws.Range("L122").Value = 23488.2
$ws.Range("M122").Value = -1979.49994
$ws.Range("N122").Value = -28388.2
$ws.Range("H135").Value = 1254281.9
$ws.Range("I135").Value = 388.6
$ws.Range("K135").Value = 3497.4
$ws.Range("M135").Value = -962.4000000000001

$ws = $wb.Worksheets("GSM")
$ws.Range("H70").Value = 4659.5244
$ws.Range("I70").Value = 4377.6943
$ws.Range("J70").Value = 5065.36
$ws.Range("K70").Value = 4377.6943
$ws.Range("L70").Value = 5065.36
$ws.Range("M70").Value = -4107.6943
$ws.Range("N70").Value = -5605.36
$ws.Range("H73").Value = 4659.5244
$ws.Range("I73").Value = 4377.6943
$ws.Range("J73").Value = 5065.36
$ws.Range("K73").Value = 4377.6943
$ws.Range("L73").Value = 5065.36
$ws.Range("M73").Value = -3441.6943
$ws.Range("N73").Value = -6937.36
$ws.Range("H80").Value = 2158.3333
$ws.Range("I80").Value = 2128.5715
$ws.Range("J80").Value = 2200
$ws.Range("K80").Value = 2128.5715
$ws.Range("L80").Value = 2200
$ws.Range("M80").Value = -1130.5715
$ws.Range("N80").Value = -4196
$ws.Range("H83").Value = 2158.3333
$ws.Range("I83").Value = 2128.5715
$ws.Range("J83").Value = 2200
$ws.Range("K83").Value = 10642.8575
$ws.Range("L83").Value = 11000
$ws.Range("M83").Value = -5650.8575
$ws.Range("N83").Value = -20984
$ws.Range("H93").Value = 19962.25
$ws.Range("J93").Value = 19962.25
$ws.Range("L93").Value = 19962.25
$ws.Range("N93").Value = -23706.25
$ws.Range("H102").Value = 4380.6665
$ws.Range("I102").Value = 5170.6665
$ws.Range("J102").Value = 3985.6667
$ws.Range("K102").Value = 5170.6665
$ws.Range("L102").Value = 3985.6667
$ws.Range("M102").Value = -3548.6665
$ws.Range("N102").Value = -7229.6667
$ws.Range("H122").Value = 2400.8684
$ws.Range("I122").Value = 1815.7084
$ws.Range("J122").Value = 3404
$ws.Range("K122").Value = 5447.1252
$ws.Range("L122").Value = 10212
$ws.Range("M122").Value = -2997.1252
$ws.Range("N122").Value = -15112

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 42750.32
$ws.Range("I7").Value = 73836
$ws.Range("J7").Value = 3186.7273
$ws.Range("K7").Value = 73836
$ws.Range("L7").Value = 3186.7273
$ws.Range("M7").Value = -73724
$ws.Range("N7").Value = -3410.7273
$ws.Range("H22").Value = 850.5
$ws.Range("I22").Value = 950
$ws.Range("J22").Value = 751
$ws.Range("K22").Value = 950
$ws.Range("L22").Value = 751
$ws.Range("M22").Value = -655
$ws.Range("N22").Value = -1341
$ws.Range("H27").Value = 850.5
$ws.Range("I27").Value = 950
$ws.Range("J27").Value = 751
$ws.Range("K27").Value = 950
$ws.Range("L27").Value = 751
$ws.Range("M27").Value = -843
$ws.Range("N27").Value = -965
$ws.Range("H40").Value = 58767.223
$ws.Range("J40").Value = 3415
$ws.Range("L40").Value = 3415
$ws.Range("N40").Value = -3687
$ws.Range("H45").Value = 12560
$ws.Range("I45").Value = 5746.6665
$ws.Range("K45").Value = 5746.6665
$ws.Range("M45").Value = -5339.6665
$ws.Range("H122").Value = 18520632
$ws.Range("I122").Value = 55556704
$ws.Range("J122").Value = 2595
$ws.Range("K122").Value = 166670112
$ws.Range("L122").Value = 7785
$ws.Range("M122").Value = -166667662
$ws.Range("N122").Value = -12685
$ws.Range("H126").Value = 42750.32
$ws.Range("I126").Value = 73836
$ws.Range("J126").Value = 3186.7273
$ws.Range("K126").Value = 221508
$ws.Range("L126").Value = 9560.1819
$ws.Range("M126").Value = -219038
$ws.Range("N126").Value = -14500.1819

$ws = $wb.Worksheets("WVR")
$ws.Range("H113").Value = 1266.909
$ws.Range("I113").Value = 1554.5
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 4663.5
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -2493.5
$ws.Range("N113").Value = -5840
$ws.Range("H124").Value = 30400
$ws.Range("J124").Value = 30400
$ws.Range("L124").Value = 30400
$ws.Range("N124").Value = -40220
$ws.Range("H126").Value = 167922.5
$ws.Range("I126").Value = 250687.5
$ws.Range("J126").Value = 2392.5
$ws.Range("K126").Value = 752062.5
$ws.Range("L126").Value = 7177.5
$ws.Range("M126").Value = -749592.5
$ws.Range("N126").Value = -12117.5
$ws.Range("H136").Value = 324594.84
$ws.Range("I136").Value = 589197.75
$ws.Range("J136").Value = 3291.2856
$ws.Range("K136").Value = 1767593.25
$ws.Range("L136").Value = 9873.856800000001
$ws.Range("M136").Value = -1765043.25
$ws.Range("N136").Value = -14973.8568
